# Updates the CryCompanywiseStockReport stock figures: quantity (F), value (G),
# rate (D/E) and sub/grand total (B) cells are corrected to their final
# reported numbers. For several item pairs the two report rows were
# effectively swapped (the Closing Stock values moved from one row to the
# other), which shows up below as a straightforward exchange of values
# between the two row's cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 202.88
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 81.98
$ws.Range("F37").Value = 30
$ws.Range("G37").Value = 799.2
$ws.Range("B40").Value = 51256.04
$ws.Range("F47").Value = 171
$ws.Range("G47").Value = 32984.19
$ws.Range("F48").Value = 58
$ws.Range("G48").Value = 2442.96
$ws.Range("F51").Value = 70
$ws.Range("G51").Value = 1591.8
$ws.Range("F57").Value = 115
$ws.Range("G57").Value = 10757.1
$ws.Range("B72").Value = 162623.38
$ws.Range("B132").Value = 65258
$ws.Range("B133").Value = 64196
$ws.Range("F160").Value = 273
$ws.Range("G160").Value = 9101.82
$ws.Range("B161").Value = 32413.25
$ws.Range("B167").Value = 57756
$ws.Range("E167").Value = 79.37
$ws.Range("F167").Value = -100
$ws.Range("G167").Value = -6644
$ws.Range("B168").Value = 64350
$ws.Range("E168").Value = 70.63
$ws.Range("F168").Value = 2
$ws.Range("G168").Value = 132.88
$ws.Range("F216").Value = 95
$ws.Range("G216").Value = 5168.95
$ws.Range("B224").Value = 61564.82
$ws.Range("B303").Value = 61610
$ws.Range("E303").Value = 122.71
$ws.Range("F303").Value = -58
$ws.Range("G303").Value = -5957.18
$ws.Range("B304").Value = 63565
$ws.Range("E304").Value = 109.19
$ws.Range("F304").Value = 60
$ws.Range("G304").Value = 6162.6
$ws.Range("B312").Value = 63531
$ws.Range("E312").Value = 152.53
$ws.Range("F312").Value = 20
$ws.Range("G312").Value = 2869.6
$ws.Range("B313").Value = 57802
$ws.Range("E313").Value = 162.71
$ws.Range("F313").Value = -79
$ws.Range("G313").Value = -11334.92
$ws.Range("F357").Value = 175
$ws.Range("G357").Value = 25312
$ws.Range("B362").Value = 68174.56
$ws.Range("F367").Value = 171
$ws.Range("G367").Value = 24040.89
$ws.Range("F368").Value = 14
$ws.Range("G368").Value = 10417.54
$ws.Range("B369").Value = 51808.74
$ws.Range("F371").Value = 3
$ws.Range("G371").Value = 165.99
$ws.Range("F372").Value = 41
$ws.Range("G372").Value = 2268.53
$ws.Range("B378").Value = 45232.25
$ws.Range("B387").Value = 47097
$ws.Range("D387").Value = 112.28
$ws.Range("E387").Value = 134.16
$ws.Range("F387").Value = 15
$ws.Range("G387").Value = 1684.2
$ws.Range("B388").Value = 58047
$ws.Range("D388").Value = 105.54
$ws.Range("E388").Value = 126.1
$ws.Range("F388").Value = 32
$ws.Range("G388").Value = 3377.28
$ws.Range("F393").Value = 338
$ws.Range("G393").Value = 32650.8
$ws.Range("B395").Value = 48386.81
$ws.Range("F402").Value = 106
$ws.Range("G402").Value = 2700.88
$ws.Range("F408").Value = 2
$ws.Range("G408").Value = 68.62
$ws.Range("F409").Value = 52
$ws.Range("G409").Value = 2108.08
$ws.Range("B423").Value = 152809.39
$ws.Range("F458").Value = 46
$ws.Range("G458").Value = 12476.12
$ws.Range("F461").Value = 29
$ws.Range("G461").Value = 6444.67
$ws.Range("F462").Value = 140
$ws.Range("G462").Value = 37781.8
$ws.Range("B464").Value = 78180.28
$ws.Range("F472").Value = 2
$ws.Range("G472").Value = 82.40000000000001
$ws.Range("B482").Value = 40906.91
$ws.Range("B485").Value = 53319
$ws.Range("E485").Value = 310.64
$ws.Range("F485").Value = -6
$ws.Range("G485").Value = -1643.52
$ws.Range("B486").Value = 64810
$ws.Range("E486").Value = 291.22
$ws.Range("F486").Value = 0
$ws.Range("G486").Value = 0
$ws.Range("B512").Value = 60022
$ws.Range("E512").Value = 37.22
$ws.Range("F512").Value = -113
$ws.Range("G512").Value = -3709.79
$ws.Range("B513").Value = 64830
$ws.Range("E513").Value = 34.9
$ws.Range("F513").Value = 83
$ws.Range("G513").Value = 2724.89
$ws.Range("F518").Value = 7
$ws.Range("G518").Value = 830.0599999999999
$ws.Range("F525").Value = 337
$ws.Range("G525").Value = 18494.56
$ws.Range("B531").Value = 104803.35
$ws.Range("F533").Value = 7
$ws.Range("G533").Value = 231.77
$ws.Range("F536").Value = 4
$ws.Range("G536").Value = 172.72
$ws.Range("F537").Value = 158
$ws.Range("G537").Value = 5231.38
$ws.Range("F538").Value = 2
$ws.Range("G538").Value = 86.36
$ws.Range("B541").Value = 16384.88
$ws.Range("F564").Value = 117
$ws.Range("G564").Value = 14256.45
$ws.Range("B567").Value = 16170.33
$ws.Range("F680").Value = 262
$ws.Range("G680").Value = 42734.82
$ws.Range("B686").Value = 43747.37
$ws.Range("B724").Value = 2082132.73
$ws.Range("B725").Value = 2082132.73
